$wb = $excel.ActiveWorkbook

# --- Update the "Status" text for the second handed-back file (row 3) everywhere
#     it is shown: Overview!B3 (zh-cn column), Overview!C3 (de-de column),
#     zh-cn!C3 and de-de!C3. The handback transform failed for this file, so
#     every cell that previously read "Ready for handoff" now reads
#     "Handback transform failed".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# --- Record the handback/handoff file name mismatch error detail for each locale ---
$wsZhCn.Range("K3").Value = "Handback file name: 3jg4mvee.kfg is different with handoff file name: 049e68d2-6839-4118-9d93-2b1938afda5a.75ccf7302be83bc575803cb3ec423670e3b3664d.zh-cn."
$wsDeDe.Range("K3").Value = "Handback file name: 3jg4mvee.kfg is different with handoff file name: 049e68d2-6839-4118-9d93-2b1938afda5a.75ccf7302be83bc575803cb3ec423670e3b3664d.de-de."
